# Add "- Steps for Making Ontology" to the "005 - Ontology Development 101"
# headline on the "video 005" slide (see commit: "add video 005 and 006").
#
# Original paragraph had two runs:
#   Run1 "005"                              (sz=6000)
#   Run2 " \u2013 Ontology Development 101" (sz=3600)
#
# Target paragraph has three runs:
#   Run1 "005"                                                (unchanged)
#   Run2 " \u2013 Ontology Development 101 \u2013 Steps for " (extended)
#   Run3 "Making Ontology"                                    (new, same look)

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt.Contains("Ontology Development 101")) {
                $targetShape = $shape
                $targetSlide = $slide
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate " - Ontology Development 101" (the run following the "005" run) by
# its position in the full text and grab it as a sub-range so we only touch
# that run, leaving the "005" run alone.
$fullText = $tr.Text
$marker = " " + [char]0x2013 + " Ontology Development 101"
$startPos = $fullText.IndexOf($marker) + 1
$run2 = $tr.Characters($startPos, $marker.Length)

$run2.Text = " " + [char]0x2013 + " Ontology Development 101 " + [char]0x2013 + " Steps for "

# Append the new trailing run; InsertAfter creates it with the same
# character formatting as the run immediately before it (sz=3600, bold,
# tx1/lumMod 75% fill) which matches the extended run above.
$newRun = $tr.InsertAfter("Making Ontology")
